{"js": "// Remove the redundant \"(\" prefix and \" hari sejak di terbitkan)\" suffix\n// around the \"${jmlhari}\" placeholder in the \"Rekomendasi berlaku\" row,\n// turning \"(${jmlhari} hari sejak di terbitkan)\" into \"${jmlhari}\".\n\nconst body = context.document.body;\n\n// Locate the paragraph that holds the \"jmlhari\" placeholder.\nconst marker = body.search(\"jmlhari\", { matchCase: true });\nmarker.load(\"items\");\nawait context.sync();\n\nif (marker.items.length === 0) {\n  throw new Error('Could not find \"jmlhari\" placeholder in the document.');\n}\n\nconst markerParagraph = marker.items[0].paragraphs.getFirst();\nconst paragraphRange = markerParagraph.getRange();\n\n// Scope the two searches to just this paragraph so we don't touch any of\n// the many other parentheses / text elsewhere in the document.\nconst leadingParen = paragraphRange.search(\"(\", { matchCase: true });\nleadingParen.load(\"items\");\nconst trailingText = paragraphRange.search(\" hari sejak di terbitkan)\", { matchCase: true });\ntrailingText.load(\"items\");\nawait context.sync();\n\nif (trailingText.items.length === 0 || leadingParen.items.length === 0) {\n  throw new Error(\"Expected surrounding text around the jmlhari placeholder was not found.\");\n}\n\n// Delete the trailing text first, then the leading \"(\" \u2014 deleting in this\n// (tail-to-head) order keeps the earlier range's offsets valid since\n// nothing before it shifts.\ntrailingText.items[0].delete();\nawait context.sync();\n\nleadingParen.items[0].delete();\nawait context.sync();\n", "ps1": "# Remove the redundant \"(\" prefix and \" hari sejak di terbitkan)\" suffix\n# around the \"${jmlhari}\" placeholder in the \"Rekomendasi berlaku\" row,\n# turning \"(${jmlhari} hari sejak di terbitkan)\" into \"${jmlhari}\".\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the \"jmlhari\" placeholder.\n$paras = $d.Paragraphs\n$count = $paras.Count\n$target = $null\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -like \"*jmlhari*\") {\n        $target = $p\n        break\n    }\n}\nif ($null -eq $target) {\n    throw \"Could not find the paragraph containing the 'jmlhari' placeholder.\"\n}\n\n$pStart = $target.Range.Start\n$pEnd = $target.Range.End\n\n# Note: document-level Ranges created via $d.Range(start, end) reliably\n# support .Delete(); Ranges read from Paragraph.Range do not, so route all\n# mutations through $d.Range(...).\n\n# Delete the trailing \" hari sejak di terbitkan)\" first - removing the tail\n# before the head keeps the head's offsets valid (nothing earlier shifts).\n$tailRange = $d.Range($pStart, $pEnd)\n$tailFound = $tailRange.Find.Execute(\" hari sejak di terbitkan)\")\nif (-not $tailFound) {\n    throw \"Could not find the trailing text ' hari sejak di terbitkan)'.\"\n}\n$tailRange.Delete()\n\n# Re-derive the (now shorter) paragraph end, then delete the leading \"(\".\n$pEnd2 = $target.Range.End\n$headRange = $d.Range($pStart, $pEnd2)\n$headFound = $headRange.Find.Execute(\"(\")\nif (-not $headFound) {\n    throw \"Could not find the leading '('.\"\n}\n$headRange.Delete()\n"}
